$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move column J (old "Kernel time (ms)" values) to column K ---
# Row 1: header cell ("loop unrolling"), keep its style (yellow fill, s=2)
$ws.Range("J1").Copy($ws.Range("K1"))
$ws.Range("J1").Clear()

# Row 2: column label "Kernel time (ms)"
$ws.Range("J2").Copy($ws.Range("K2"))
$ws.Range("J2").Clear()

# Data rows: J -> K (skip row 7, which has no data)
$dataRows = @(3, 4, 5, 6, 8, 9, 10)
foreach ($r in $dataRows) {
    $ws.Range("J$r").Copy($ws.Range("K$r"))
    $ws.Range("J$r").Clear()
}

# --- Add new "GOP/s" columns I and L ---
# Column I: GOP/s computed from column H (Kernel time (ms))
# Column L: GOP/s computed from column K (Kernel time (ms), multi-thread init)
$ws.Range("I2").Value = "GOP/s"
$ws.Range("L2").Value = "GOP/s"

$ws.Range("I3").Formula = "=16*16*128*128*4*4*2/H3*1000/1000000000"
$ws.Range("L3").Formula = "=16*16*128*128*4*4*2/K3*1000/1000000000"

$ws.Range("I4:I10").Formula = "=16*16*128*128*4*4*2/H4*1000/1000000000"
$ws.Range("L4:L10").Formula = "=16*16*128*128*4*4*2/K4*1000/1000000000"

# Remove stray cells auto-created in the empty row 7 (no source data there)
$ws.Range("I7").Clear()
$ws.Range("L7").Clear()

# --- Sheet view adjustments ---
$ws.Range("N7").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 4

# --- Column widths for new columns ---
$ws.Columns("I").ColumnWidth = 13.83203125
$ws.Columns("K").ColumnWidth = 14.08203125
$ws.Columns("L").ColumnWidth = 14
